# Adds a "Github link: https://github.com/nazzy243/PRT582-assignment-1.git"
# sentence to the (previously empty) paragraph that carries the "_GoBack"
# bookmark, matching the upstream commit's OOXML diff exactly:
#   - a <w:proofErr type="spellStart"/> ... <w:proofErr type="spellEnd"/>
#     pair wrapping a run that reads "Github"
#   - a run with " link: "
#   - a run with the bare GitHub URL
#   - all three runs use the same Times New Roman / 24 half-point run
#     properties already used throughout the document
#   - the existing <w:bookmarkStart w:name="_GoBack"/><w:bookmarkEnd/> pair
#     is left in place, immediately after the new runs

$d = $word.ActiveDocument

# Locate the target paragraph robustly via the "_GoBack" bookmark (rather
# than a hard-coded paragraph index) and expand to its full paragraph range.
$bm = $d.Bookmarks("_GoBack")
$target = $bm.Range.Duplicate
[void]$target.Expand(4)  # wdParagraph

# Rebuild that single paragraph in place (same w14:paraId / rsid* / pPr as
# before) with the new proofErr-wrapped runs inserted ahead of the bookmark.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        '<w:body>' + `
          '<w:p w14:paraId="225D9C34" w14:textId="7881A318" w:rsidR="004B2910" w:rsidRPr="008A47F7" w:rsidRDefault="004B2910" w:rsidP="004B2910">' + `
            '<w:pPr>' + `
              '<w:spacing w:after="0" w:line="480" w:lineRule="auto"/>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
            '</w:pPr>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t>Github</w:t>' + `
            '</w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t xml:space="preserve"> link: </w:t>' + `
            '</w:r>' + `
            '<w:r>' + `
              '<w:rPr>' + `
                '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
                '<w:sz w:val="24"/>' + `
                '<w:szCs w:val="24"/>' + `
              '</w:rPr>' + `
              '<w:t>https://github.com/nazzy243/PRT582-assignment-1.git</w:t>' + `
            '</w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
            '<w:bookmarkEnd w:id="0"/>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

[void]$target.InsertXML($xml)
